$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 0.768783636774762
$ws.Range("K2").Value = 0.0322867007401506
$ws.Range("L2").Value = -0.0905953164835239
$ws.Range("M2").Value = 0.06859142398971679
$ws.Range("N2").Value = 7.87480505857331
$ws.Range("P2").Value = "Likely improving"

# Row 3
$ws.Range("F3").Value = 0.59675202974633
$ws.Range("K3").Value = 0.444469079618711
$ws.Range("L3").Value = -5.23510506637395
$ws.Range("M3").Value = 4.64307616797114
$ws.Range("N3").Value = 0.399702409729057
$ws.Range("P3").Value = "As likely as not improving"

# Row 4
$ws.Range("F4").Value = 0.768783636774762
$ws.Range("K4").Value = 0.362915400033441
$ws.Range("L4").Value = -0.830349856934493
$ws.Range("M4").Value = 0.957115847805783
$ws.Range("N4").Value = 7.25830800066882
$ws.Range("P4").Value = "Likely improving"
